# Append four new match rows (rows 3-6) to the Deepak Hooda stats sheet.
# Numeric-looking stat columns (G:K) are written with a leading
# quote-prefix so Excel stores them as TEXT (matching the existing
# sheet convention where every stat is a text value), instead of
# silently re-typing them as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("A3").Value = ' Oct 30 2020'
$ws.Range("B3").Value = ' Abu Dhabi'
$ws.Range("C3").Value = 'Royals won by 7 wickets (with 15 balls remaining)'
$ws.Range("D3").Value = 'Kings XI Punjab'
$ws.Range("E3").Value = 'Rajasthan Royals'
$ws.Range("F3").Value = 'Deepak Hooda '
$ws.Range("G3").Value = '''1'
$ws.Range("H3").Value = '''1'
$ws.Range("I3").Value = '''0'
$ws.Range("J3").Value = '''0'
$ws.Range("K3").Value = '''100.00'

# Row 4
$ws.Range("A4").Value = ' Oct 20 2020'
$ws.Range("B4").Value = ' Dubai (DSC)'
$ws.Range("C4").Value = 'Kings XI won by 5 wickets (with 6 balls remaining)'
$ws.Range("D4").Value = 'Kings XI Punjab'
$ws.Range("E4").Value = 'Delhi Capitals'
$ws.Range("F4").Value = 'Deepak Hooda '
$ws.Range("G4").Value = '''15'
$ws.Range("H4").Value = '''22'
$ws.Range("I4").Value = '''1'
$ws.Range("J4").Value = '''0'
$ws.Range("K4").Value = '''68.18'

# Row 5
$ws.Range("A5").Value = ' Oct 18 2020'
$ws.Range("B5").Value = ' Dubai (DSC)'
$ws.Range("C5").Value = 'Match tied (Kings XI won the one-over eliminator)'
$ws.Range("D5").Value = 'Kings XI Punjab'
$ws.Range("E5").Value = 'Mumbai Indians'
$ws.Range("F5").Value = 'Deepak Hooda '
$ws.Range("G5").Value = '''23'
$ws.Range("H5").Value = '''16'
$ws.Range("I5").Value = '''1'
$ws.Range("J5").Value = '''1'
$ws.Range("K5").Value = '''143.75'

# Row 6
$ws.Range("A6").Value = ' Oct 24 2020'
$ws.Range("B6").Value = ' Dubai (DSC)'
$ws.Range("C6").Value = 'Kings XI won by 12 runs'
$ws.Range("D6").Value = 'Kings XI Punjab'
$ws.Range("E6").Value = 'Sunrisers Hyderabad'
$ws.Range("F6").Value = 'Deepak Hooda '
$ws.Range("G6").Value = '''0'
$ws.Range("H6").Value = '''2'
$ws.Range("I6").Value = '''0'
$ws.Range("J6").Value = '''0'
$ws.Range("K6").Value = '''0.00'
